# Apply the "Updating jacoco csv and Codec Calculations (All Metrics)" edit:
# Fill in the previously-empty coverage / churn metric values for the
# 1.12-RC1 / RC2 / RC3 columns (D:F) in both the "Coverage Results"
# table (rows 6-8) and the "Churned Code Results" table (rows 13-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Codec")

# --- Coverage Results (Metric 1, 2 & 4) block ---
# Code Coverage(%)
$ws.Range("D6").Value = 96
$ws.Range("E6").Value = 96
$ws.Range("F6").Value = 96

# Branch Coverage(%)
$ws.Range("D7").Value = 92
$ws.Range("E7").Value = 91
$ws.Range("F7").Value = 89

# Complexity
$ws.Range("D8").Value = 1733
$ws.Range("E8").Value = 1868
$ws.Range("F8").Value = 1989

# --- Churned Code Results (Metric 5) block ---
# No. Of Classes
$ws.Range("D13").Value = 45
$ws.Range("E13").Value = 50
$ws.Range("F13").Value = 52

# Line Coverage(%)
$ws.Range("D14").Value = 95
$ws.Range("E14").Value = 92
$ws.Range("F14").Value = 92

# Mutation Coverage(%)
$ws.Range("D15").Value = 90
$ws.Range("E15").Value = 87
$ws.Range("F15").Value = 87

# Restore the active selection to match the saved workbook state
$ws.Range("D15").Select()
